$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "ACL_TD11 / Datta" record (previously row 11) and the "ACL_TD2 / Ashutosh"
# record (previously row 2) swap places. Rows 3-10 keep their person data but
# get a freshly-numbered WebSite value in column G.

$ws.Range("B2").Value = "Datta"
$ws.Range("C2").Value = "centera"
$ws.Range("D2").Value = "datta@detore.co"
$ws.Range("E2").Value = 1992342396
$ws.Range("F2").Value = "11-1234576"
$ws.Range("G2").Value = "temp.clidiem10.com"

$ws.Range("G3").Value = "temp.clidiem2.com"
$ws.Range("G4").Value = "temp.clidiem3.com"
$ws.Range("G5").Value = "temp.clidiem4.com"
$ws.Range("G6").Value = "temp.clidiem5.com"
$ws.Range("G7").Value = "temp.clidiem6.com"
$ws.Range("G8").Value = "temp.clidiem7.com"
$ws.Range("G9").Value = "temp.clidiem8.com"
$ws.Range("G10").Value = "temp.clidiem9.com"

$ws.Range("B11").Value = "Ashutosh"
$ws.Range("C11").Value = "Panacy"
$ws.Range("D11").Value = "ashutosh@glsie.co"
$ws.Range("E11").Value = 1992342387
$ws.Range("F11").Value = "11-1234567"
$ws.Range("G11").Value = "temp.clidiem1.com"

# Rebuild the hyperlinks so each cell points at the record that now lives
# there (D/H/I hold mailto links; H/I always point at the shared "Test@123"
# placeholder account).
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("D11"), "mailto:ashutosh@glsie.co")
$ws.Hyperlinks.Add($ws.Range("H11"), "mailto:Test@123")
$ws.Hyperlinks.Add($ws.Range("I11"), "mailto:Test@123")
$ws.Hyperlinks.Add($ws.Range("H3"), "mailto:Test@123")
$ws.Hyperlinks.Add($ws.Range("H4"), "mailto:Test@123")
$ws.Hyperlinks.Add($ws.Range("H5"), "mailto:Test@123")
$ws.Hyperlinks.Add($ws.Range("H6"), "mailto:Test@123")
$ws.Hyperlinks.Add($ws.Range("H7"), "mailto:Test@123")
$ws.Hyperlinks.Add($ws.Range("H8"), "mailto:Test@123")
$ws.Hyperlinks.Add($ws.Range("H9"), "mailto:Test@123")
$ws.Hyperlinks.Add($ws.Range("H10"), "mailto:Test@123")
$ws.Hyperlinks.Add($ws.Range("H2"), "mailto:Test@123")
$ws.Hyperlinks.Add($ws.Range("I3"), "mailto:Test@123")
$ws.Hyperlinks.Add($ws.Range("I4"), "mailto:Test@123")
$ws.Hyperlinks.Add($ws.Range("I5"), "mailto:Test@123")
$ws.Hyperlinks.Add($ws.Range("I6"), "mailto:Test@123")
$ws.Hyperlinks.Add($ws.Range("I7"), "mailto:Test@123")
$ws.Hyperlinks.Add($ws.Range("I8"), "mailto:Test@123")
$ws.Hyperlinks.Add($ws.Range("I9"), "mailto:Test@123")
$ws.Hyperlinks.Add($ws.Range("I10"), "mailto:Test@123")
$ws.Hyperlinks.Add($ws.Range("I2"), "mailto:Test@123")
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:amet@sedt.edu")
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:amron@seiwr.co")
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:koyal@12coms.co")
$ws.Hyperlinks.Add($ws.Range("D7"), "mailto:lasa32@cosms.co")
$ws.Hyperlinks.Add($ws.Range("D8"), "mailto:tina@nullasw.co")
$ws.Hyperlinks.Add($ws.Range("D9"), "mailto:toshirew@agetew.co")
$ws.Hyperlinks.Add($ws.Range("D10"), "mailto:desh@teshart.co")
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:datta@detore.co")
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:velit@nams.com")

# Restore the sheet selection/view state recorded after the edit.
$ws.Range("B2:I2").Select()
